$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before "Thickness" (old column F) to hold the
# referenced configuration name, shifting Thickness..File Name right by one.
$ws.Columns("F").Insert()

# New header cell: "Configuration"
$ws.Range("F1").Value = "Configuration"

# Row 1 formatting: taller header row to fit the extra wrapped column text
$ws.Rows("1").RowHeight = 44.25

# Column width for the freshly inserted "Configuration" column
$ws.Columns("F").ColumnWidth = 17

# Nudge the other column widths to the refreshed layout
$ws.Columns("A").ColumnWidth = 8.6666666666667
$ws.Columns("B").ColumnWidth = 25.6666666666667
$ws.Columns("C").ColumnWidth = 7
$ws.Columns("D").ColumnWidth = 44.1666666666667
$ws.Columns("E").ColumnWidth = 39.6666666666667
$ws.Columns("G").ColumnWidth = 9.8333333333333
$ws.Columns("H").ColumnWidth = 22.8333333333333
$ws.Columns("I").ColumnWidth = 9.1666666666667
$ws.Columns("J").ColumnWidth = 12.6666666666667

# Re-point the hidden _FilterDatabase defined name at the widened header row
$fd = $wb.Names.Item("Parts!_FilterDatabase")
$fd.RefersTo = "=Parts!`$A`$1:`$J`$1"

# Resize the AutoFilter to the new last column
$ws.Range("A1:J1").AutoFilter()

# Selection left on D1, matching the edited workbook's last active cell
$ws.Range("D1").Select()
